# correct rotation of p&p'd parts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# LED1 (row 75): was "Green 1.85V" / "CHIP-LED0805" -> "white led 0805 C34499" / "LED 0805"
$ws.Range("A75").Value = "white led 0805 C34499"
$ws.Range("C75").Value = "LED 0805"

# LED2 (row 76): was "Orange 1.85V" / "CHIP-LED0805" -> "yellow led 0805 C2296" / "LED 0805"
$ws.Range("A76").Value = "yellow led 0805 C2296"
$ws.Range("C76").Value = "LED 0805"

# LED3 (row 77): was "Red 1.85V" / "CHIP-LED0805" -> "red led 0805 C84256" / "LED 0805"
$ws.Range("A77").Value = "red led 0805 C84256"
$ws.Range("C77").Value = "LED 0805"

# LED4 (row 78): was "Green 1.85V" / "CHIP-LED0805" -> "white led 0805 C34499" / "LED 0805"
$ws.Range("A78").Value = "white led 0805 C34499"
$ws.Range("C78").Value = "LED 0805"

# Row 154 (S2 / KSS221GLFS / KSS) is removed entirely - the whole part was dropped from the BoM.
$ws.Rows("154").Delete()

# Update the view state to match the author's final selection/scroll position.
$ws.Range("B77").Select()
$excel.ActiveWindow.ScrollRow = 64
